$wb = $excel.ActiveWorkbook

$updates = @{
    "hybNew.tot" = @{
        "C3" = 267.48369273734477
        "D3" = 206.6269245827927
        "E3" = 126.4995755972538
        "F3" = 0.04416252474882043
        "G3" = 73.4482571093704
        "J3" = 140.98411714009097
        "L3" = 1.114502688838678
        "M3" = 0.6334199036417827
    }
    "hybNew.new" = @{
        "C3" = 23.04354067941053
        "D3" = 19.05279540334587
        "E3" = 13.363421340464242
        "F3" = 6.2949378084119605
        "G3" = 11.417344210771299
        "I3" = 12.757857594933908
        "J3" = 9.680119338946287
        "K3" = 5.689374062881626
        "L3" = 0.7243743269274185
        "M3" = 0.42574232435927817
    }
    "hybAll.tot" = @{
        "C3" = 154.0160415479052
        "D3" = 122.23226258722894
        "E3" = 80.82722532025296
        "F3" = 12.516088600210448
        "G3" = 52.63654792747733
        "J3" = 73.18881622765223
        "K3" = 41.40503726697598
        "L3" = 0.9054970764821395
        "M3" = 0.5122659735370265
    }
    "hybNew.tot.SR0" = @{
        "C3" = 273.87934009590884
        "D3" = 211.0302574539941
        "E3" = 128.01321226685042
        "F3" = 1.4180575439565928
        "G3" = 76.50242746632463
        "K3" = 83.01704518714368
        "L3" = 1.1394615074965282
        "M3" = 0.6485037264285674
    }
    "hybNew.new.SR0" = @{
        "C3" = 27.698766812872968
        "D3" = 22.6365444358682
        "E3" = 16.034972912737253
        "F3" = 7.90382758796596
        "G3" = 14.122190961836736
        "H3" = 12.05127931820565
        "I3" = 14.732716847902239
        "J3" = 11.663793900135715
        "L3" = 0.7273971689014
        "M3" = 0.4116983270914689
    }
    "hybAll.tot.SR0" = @{
        "C3" = 161.17179097640673
        "D3" = 127.0639293146295
        "E3" = 82.16595190452469
        "F3" = 13.126269968120003
        "G3" = 55.697206454464386
        "I3" = 113.9376593465095
        "J3" = 79.00583907188204
        "K3" = 44.897977410104815
        "L3" = 0.9615398743713863
        "M3" = 0.5464304419216786
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cells = $updates[$sheetName]
    foreach ($cellRef in $cells.Keys) {
        $ws.Range($cellRef).Value = $cells[$cellRef]
    }
}
